# Update numeric results in Sheet1 (KNN imputation result data) to new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -10.96
$ws.Range("D3").Value = -7.348000000000001
$ws.Range("E8").Value = 16.606
$ws.Range("E11").Value = 17.409
$ws.Range("B12").Value = 4.872
$ws.Range("C14").Value = -12.987
$ws.Range("E14").Value = 16.783
$ws.Range("E15").Value = 15.916
$ws.Range("C26").Value = -12.427
$ws.Range("D30").Value = -7.25
$ws.Range("C31").Value = -12.772
$ws.Range("B32").Value = 6.048999999999999
$ws.Range("C35").Value = -12.311
$ws.Range("B36").Value = 8.821
$ws.Range("E36").Value = 16.396
$ws.Range("C37").Value = -13.387
$ws.Range("B38").Value = 5.203
$ws.Range("D44").Value = -7.409999999999999
$ws.Range("C45").Value = -12.67
$ws.Range("B46").Value = 5.553
$ws.Range("B54").Value = 5.275
$ws.Range("B55").Value = 4.599000000000001
$ws.Range("C57").Value = -13.697
$ws.Range("D58").Value = -8.035
$ws.Range("E64").Value = 17.444
$ws.Range("B67").Value = 5.327999999999999
$ws.Range("B69").Value = 5.1
$ws.Range("B72").Value = 5.262
$ws.Range("D84").Value = -8.238
$ws.Range("D89").Value = -7.363
$ws.Range("E89").Value = 17.235
$ws.Range("B91").Value = 5.327
$ws.Range("D91").Value = -6.635000000000001
$ws.Range("D92").Value = -6.663000000000001
$ws.Range("B99").Value = 5.718
$ws.Range("C100").Value = -12.662
$ws.Range("C102").Value = -13.583
$ws.Range("D102").Value = -7.723999999999999
